$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2042682926829268
$ws.Range("C2").Value = 0.5548780487804879
$ws.Range("J2").Value = 0.01829268292682927
$ws.Range("P2").Value = 0.1432926829268293
$ws.Range("S2").Value = 0.07926829268292683
$ws.Range("B3").Value = 0.01036269430051814
$ws.Range("C3").Value = 0.02590673575129534
$ws.Range("J3").Value = 0.0310880829015544
$ws.Range("P3").Value = 0.7668393782383419
$ws.Range("S3").Value = 0.1658031088082902
$ws.Range("J4").Value = 0.06382978723404255
$ws.Range("P4").Value = 0.6808510638297872
$ws.Range("S4").Value = 0.2553191489361702
$ws.Range("B6").Value = 0.08095238095238096
$ws.Range("F6").Value = 0.05714285714285714
$ws.Range("J6").Value = 0.2714285714285714
$ws.Range("O6").Value = 0.01428571428571429
$ws.Range("Q6").Value = 0.1666666666666667
$ws.Range("R6").Value = 0.07142857142857142
$ws.Range("S6").Value = 0.3380952380952381
$ws.Range("B7").Value = 0.1319444444444444
$ws.Range("D7").Value = 0.02777777777777778
$ws.Range("F7").Value = 0.05555555555555555
$ws.Range("J7").Value = 0.09027777777777778
$ws.Range("O7").Value = 0.02083333333333333
$ws.Range("Q7").Value = 0.2430555555555556
$ws.Range("R7").Value = 0.06944444444444445
$ws.Range("S7").Value = 0.3611111111111111
$ws.Range("B8").Value = 0.1101321585903084
$ws.Range("D8").Value = 0.01762114537444934
$ws.Range("F8").Value = 0.05286343612334802
$ws.Range("J8").Value = 0.1123348017621145
$ws.Range("O8").Value = 0.01762114537444934
$ws.Range("Q8").Value = 0.2048458149779736
$ws.Range("R8").Value = 0.1079295154185022
$ws.Range("S8").Value = 0.3766519823788546
$ws.Range("B9").Value = 0.07567567567567568
$ws.Range("D9").Value = 0.02162162162162162
$ws.Range("F9").Value = 0.06486486486486487
$ws.Range("J9").Value = 0.1297297297297297
$ws.Range("O9").Value = 0.01621621621621622
$ws.Range("Q9").Value = 0.1675675675675676
$ws.Range("R9").Value = 0.1189189189189189
$ws.Range("S9").Value = 0.4054054054054054
$ws.Range("B10").Value = 0.1181102362204724
$ws.Range("D10").Value = 0.02440944881889764
$ws.Range("F10").Value = 0.05984251968503937
$ws.Range("J10").Value = 0.1125984251968504
$ws.Range("O10").Value = 0.01811023622047244
$ws.Range("Q10").Value = 0.2409448818897638
$ws.Range("R10").Value = 0.1102362204724409
$ws.Range("S10").Value = 0.315748031496063
$ws.Range("G11").Value = 0.1254901960784314
$ws.Range("J11").Value = 0.1137254901960784
$ws.Range("K11").Value = 0.1882352941176471
$ws.Range("L11").Value = 0.5568627450980392
$ws.Range("S11").Value = 0.01568627450980392
$ws.Range("G12").Value = 0.7222222222222222
$ws.Range("J12").Value = 0.1944444444444444
$ws.Range("K12").Value = 0.006944444444444444
$ws.Range("L12").Value = 0.03472222222222222
$ws.Range("S12").Value = 0.04166666666666666
$ws.Range("G13").Value = 0.4545454545454545
$ws.Range("J13").Value = 0.4242424242424243
$ws.Range("S13").Value = 0.1212121212121212
$ws.Range("F15").Value = 0.03896103896103896
$ws.Range("H15").Value = 0.1645021645021645
$ws.Range("I15").Value = 0.06060606060606061
$ws.Range("J15").Value = 0.3852813852813853
$ws.Range("K15").Value = 0.04761904761904762
$ws.Range("M15").Value = 0.004329004329004329
$ws.Range("O15").Value = 0.06493506493506493
$ws.Range("S15").Value = 0.2337662337662338
$ws.Range("F16").Value = 0.01376146788990826
$ws.Range("H16").Value = 0.2064220183486239
$ws.Range("I16").Value = 0.08256880733944955
$ws.Range("J16").Value = 0.3486238532110092
$ws.Range("K16").Value = 0.1009174311926606
$ws.Range("M16").Value = 0.009174311926605505
$ws.Range("O16").Value = 0.0963302752293578
$ws.Range("S16").Value = 0.1422018348623853
$ws.Range("F17").Value = 0.0199203187250996
$ws.Range("H17").Value = 0.203187250996016
$ws.Range("I17").Value = 0.1095617529880478
$ws.Range("J17").Value = 0.3924302788844621
$ws.Range("K17").Value = 0.08565737051792828
$ws.Range("M17").Value = 0.01593625498007968
$ws.Range("O17").Value = 0.06772908366533864
$ws.Range("S17").Value = 0.1055776892430279
$ws.Range("F18").Value = 0.008620689655172414
$ws.Range("H18").Value = 0.2413793103448276
$ws.Range("I18").Value = 0.04741379310344827
$ws.Range("J18").Value = 0.4396551724137931
$ws.Range("K18").Value = 0.06465517241379311
$ws.Range("M18").Value = 0.01293103448275862
$ws.Range("O18").Value = 0.08189655172413793
$ws.Range("S18").Value = 0.103448275862069
$ws.Range("F19").Value = 0.0225498699045967
$ws.Range("H19").Value = 0.186470078057242
$ws.Range("I19").Value = 0.07718993928881179
$ws.Range("J19").Value = 0.385949696444059
$ws.Range("K19").Value = 0.0997398091934085
$ws.Range("M19").Value = 0.01734605377276669
$ws.Range("O19").Value = 0.06764960971379011
$ws.Range("S19").Value = 0.1431049436253252
